$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 9: "The Lord of the Rings Sketchbook" (JUdOAAAACAAJ) ---
$ws.Rows(9).Insert()
$ws.Range("A9").Value = 'JUdOAAAACAAJ'
$ws.Range("B9").Value = 'The Lord of the Rings Sketchbook'
$ws.Range("C9").Value = '"In The Lord of the Rings Sketchbook Alan Lee reveals in pictures and in words how he created the watercolor paintings for the special centenary edition of The Lord of the Rings. These images would prove so powerful and evocative that they would eventually define the look of Peter Jackson''s movie trilogy and earn him a coveted Academy Award." "The book is filled with more than 150 of his sketches and early conceptual pieces showing how the project progressed from idea to finished art. It also contains a selection of full-page paintings reproduced in full color, together with numerous examples of previously unseen conceptual art produced for the films and many new works drawn specially for this book." "The Lord of the Rings Sketchbook provides an insight into the imagination of the man who painted Tolkien''s vision, first on the page and then in three dimensions on the movie screen. It will also be of interest to many of the thousands of people who have bought the illustrated Lord of the Rings as well as to budding artists who want to unlock the secrets of book illustration."--BOOK JACKET.'
# Force column D to text so the numeric-looking year is not reinterpreted, then restore default styling
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2005'
$ws.Range("D9").Style = $ws.Range("D8").Style
$ws.Range("E9").Value = 'Alan Lee'

# The previous row 16 book (WBOxAQAACAAJ) is no longer present in the export; it is now at row 17 after the insert above
$ws.Rows(17).Delete()

# --- Insert a new row at position 19: "The Fellowship Of The Ring" (aWZzLPhY4o0C) ---
$ws.Rows(19).Insert()
$ws.Range("A19").Value = 'aWZzLPhY4o0C'
$ws.Range("B19").Value = 'The Fellowship Of The Ring'
$ws.Range("C19").Value = 'Begin your journey into Middle-earth... The inspiration for the upcoming original series on Prime Video, The Lord of the Rings: The Rings of Power. The Fellowship of the Ring is the first part of J.R.R. Tolkien’s epic adventure The Lord of the Rings. One Ring to rule them all, One Ring to find them, One Ring to bring them all and in the darkness bind them. Sauron, the Dark Lord, has gathered to him all the Rings of Power—the means by which he intends to rule Middle-earth. All he lacks in his plans for dominion is the One Ring—the ring that rules them all—which has fallen into the hands of the hobbit, Bilbo Baggins. In a sleepy village in the Shire, young Frodo Baggins finds himself faced with an immense task, as his elderly cousin Bilbo entrusts the Ring to his care. Frodo must leave his home and make a perilous journey across Middle-earth to the Cracks of Doom, there to destroy the Ring and foil the Dark Lord in his evil purpose.'
# Force column D to text so the date-looking publishedDate is not reinterpreted as a date serial, then restore default styling
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2012-02-15'
$ws.Range("D19").Style = $ws.Range("D18").Style
$ws.Range("E19").Value = 'J.R.R. Tolkien'

# The previous row 22 book (eqPUjwEACAAJ) is no longer present in the export; it is now at row 23 after the insert above
$ws.Rows(23).Delete()
